# Adding the changes we made on may 9th
#
# Prepend 13 new gyroscope readings above the existing data (they become
# the new rows 2..14), pushing the prior readings down, and drop the
# oldest 3 readings that fall off the end of the tracked window so the
# sheet keeps a fixed-size rolling window of 30 data rows (plus header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(0.0594066455960273, 0.1411098688840866, 0.0326812900602817),
    @(0.0723875313997268, 0.0080939643085002, 0.0835358202457428),
    @(0.0029016099870204, 0.0274889357388019, 0.0937678143382072),
    @(0.0433714315295219, 0.0073303831741213, -0.0807869285345077),
    @(0.0352774672210216, 0.0056505035609006, 0.0056505035609006),
    @(0.0557414554059505, 0.0113010071218013, -0.1820378452539444),
    @(-0.0462730415165424, -0.0134390350431203, -0.0310014113783836),
    @(0.0276416521519422, 0.001527163083665, -0.0335975885391235),
    @(-0.0464257597923278, -0.0105374250560998, -0.06368270516395561),
    @(-0.0219911485910415, -0.0183259565383195, 0.0233655963093042),
    @(-0.0684169083833694, -0.0335975885391235, 0.0587957799434661),
    @(0.0198531206697225, -0.0583376325666904, 0.0022907445672899),
    @(-0.0027488935738801, -0.0503963828086853, 0.0137444678694009)
)
$insertCount = $newData.Count

# Remember the current row 2 (it will slide down below the new block).
$oldA2 = $ws.Cells.Item(2, 1).Value2
$oldB2 = $ws.Cells.Item(2, 2).Value2
$oldC2 = $ws.Cells.Item(2, 3).Value2

# Insert the new rows starting at row 3 rather than row 2: inserting
# immediately below the styled header row causes the fresh rows to
# inherit the header's bold/centered formatting, which the original data
# rows don't have. Leaving row 2 untouched as a buffer avoids that, and
# we copy its old contents back into place afterwards.
$lastInsertRow = 2 + $insertCount
$ws.Range("A3:A$lastInsertRow").EntireRow.Insert()

# The old row 2 reading now belongs right after the new block.
$destRow = $insertCount + 2
$ws.Cells.Item($destRow, 1).Value = $oldA2
$ws.Cells.Item($destRow, 2).Value = $oldB2
$ws.Cells.Item($destRow, 3).Value = $oldC2

# Write the new readings into rows 2..(insertCount + 1), overwriting the
# old row 2 value along the way.
$r = 2
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Drop the oldest 3 readings that now sit past the tracked window (the
# original rows 19..21, which have shifted down to rows 32..34).
$ws.Range("A32:A34").EntireRow.Delete()
